# Generate Report for Handoff
# Adds two new tracked files (38b97b2f-... and 8ee410fb-...) to the
# localization-status workbook: one new row per file on the "Overview"
# sheet, and one new row per file on each of the "zh-cn" / "de-de" sheets.

$wb = $excel.ActiveWorkbook

$linkColor = 15570276  # BGR int for RGB FF6495ED (matches existing HyperLink style)

function Style-AsLink($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $linkColor
}

function Style-AsDate($range) {
    $range.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 4: 38b97b2f-13f0-44da-b781-6dbec8b0f516.md
$wsOverview.Range("A4").Value = "38b97b2f-13f0-44da-b781-6dbec8b0f516.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cdecf1e90d7fcae761894fd034b670f0906c2201/e2e/38b97b2f-13f0-44da-b781-6dbec8b0f516.md", "", "", "38b97b2f-13f0-44da-b781-6dbec8b0f516.md")
Style-AsLink($wsOverview.Range("A4"))
$wsOverview.Range("B4").Value = "Ready for handoff"
$wsOverview.Range("C4").Value = "Ready for handoff"
$wsOverview.Range("D4").Value = "2016-03-25 07:47:41"
Style-AsDate($wsOverview.Range("D4"))

# Row 5: 8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md
$wsOverview.Range("A5").Value = "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cdecf1e90d7fcae761894fd034b670f0906c2201/e2e/8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md", "", "", "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md")
Style-AsLink($wsOverview.Range("A5"))
$wsOverview.Range("B5").Value = "Ready for handoff"
$wsOverview.Range("C5").Value = "Ready for handoff"
$wsOverview.Range("D5").Value = "2016-03-25 07:47:41"
Style-AsDate($wsOverview.Range("D5"))

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Reference Tokens |
#   Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 4: 38b97b2f-13f0-44da-b781-6dbec8b0f516
$wsZhCn.Range("A4").Value = "38b97b2f-13f0-44da-b781-6dbec8b0f516.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cdecf1e90d7fcae761894fd034b670f0906c2201/e2e/38b97b2f-13f0-44da-b781-6dbec8b0f516.md", "", "", "38b97b2f-13f0-44da-b781-6dbec8b0f516.md")
Style-AsLink($wsZhCn.Range("A4"))
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Ready for handoff"
$wsZhCn.Range("D4").Value = "38b97b2f-13f0-44da-b781-6dbec8b0f516.44178893c7a354429dbf8fc9a2489f72603f388d.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1b1433220495afbb967ef936e3ba1c67a9bd780/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/38b97b2f-13f0-44da-b781-6dbec8b0f516.44178893c7a354429dbf8fc9a2489f72603f388d.zh-cn.xlf", "", "", "38b97b2f-13f0-44da-b781-6dbec8b0f516.44178893c7a354429dbf8fc9a2489f72603f388d.zh-cn.xlf")
Style-AsLink($wsZhCn.Range("D4"))
$wsZhCn.Range("E4").Value = "2016-03-25 07:47:31"
$wsZhCn.Range("H4").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J4").Value = "Include"

# Row 5: 8ee410fb-668b-48fb-a8b8-6fd9894dd5ac
$wsZhCn.Range("A5").Value = "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cdecf1e90d7fcae761894fd034b670f0906c2201/e2e/8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md", "", "", "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md")
Style-AsLink($wsZhCn.Range("A5"))
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = "Ready for handoff"
$wsZhCn.Range("D5").Value = "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.87ba0c7937f09c893f59073b767e9834435e914e.zh-cn.xlf"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d1b1433220495afbb967ef936e3ba1c67a9bd780/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/xinjiang/ht/8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.87ba0c7937f09c893f59073b767e9834435e914e.zh-cn.xlf", "", "", "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.87ba0c7937f09c893f59073b767e9834435e914e.zh-cn.xlf")
Style-AsLink($wsZhCn.Range("D5"))
$wsZhCn.Range("E5").Value = "2016-03-25 07:47:31"
$wsZhCn.Range("H5").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("J5").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de": same column layout as "zh-cn"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 4: 38b97b2f-13f0-44da-b781-6dbec8b0f516
$wsDeDe.Range("A4").Value = "38b97b2f-13f0-44da-b781-6dbec8b0f516.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/cdecf1e90d7fcae761894fd034b670f0906c2201/e2e/38b97b2f-13f0-44da-b781-6dbec8b0f516.md", "", "", "38b97b2f-13f0-44da-b781-6dbec8b0f516.md")
Style-AsLink($wsDeDe.Range("A4"))
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Ready for handoff"
$wsDeDe.Range("D4").Value = "38b97b2f-13f0-44da-b781-6dbec8b0f516.44178893c7a354429dbf8fc9a2489f72603f388d.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e74bad424b88b11e853bad2582f8b174310add2a/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/38b97b2f-13f0-44da-b781-6dbec8b0f516.44178893c7a354429dbf8fc9a2489f72603f388d.de-de.xlf", "", "", "38b97b2f-13f0-44da-b781-6dbec8b0f516.44178893c7a354429dbf8fc9a2489f72603f388d.de-de.xlf")
Style-AsLink($wsDeDe.Range("D4"))
$wsDeDe.Range("E4").Value = "2016-03-25 07:47:41"
$wsDeDe.Range("H4").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J4").Value = "Include"

# Row 5: 8ee410fb-668b-48fb-a8b8-6fd9894dd5ac
$wsDeDe.Range("A5").Value = "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/cdecf1e90d7fcae761894fd034b670f0906c2201/e2e/8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md", "", "", "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.md")
Style-AsLink($wsDeDe.Range("A5"))
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = "Ready for handoff"
$wsDeDe.Range("D5").Value = "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.87ba0c7937f09c893f59073b767e9834435e914e.de-de.xlf"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e74bad424b88b11e853bad2582f8b174310add2a/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/xinjiang/ht/8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.87ba0c7937f09c893f59073b767e9834435e914e.de-de.xlf", "", "", "8ee410fb-668b-48fb-a8b8-6fd9894dd5ac.87ba0c7937f09c893f59073b767e9834435e914e.de-de.xlf")
Style-AsLink($wsDeDe.Range("D5"))
$wsDeDe.Range("E5").Value = "2016-03-25 07:47:41"
$wsDeDe.Range("H5").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("J5").Value = "Include"

Write-Output "Report generated for handoff: 2 new files across Overview/zh-cn/de-de sheets"
